$d = $word.ActiveDocument

# --- Locate the "Low Coupling" paragraph -----------------------------------
# Search on unique text from the paragraph we need to rewrite, then expand
# the found range to the whole paragraph (wdParagraph = 4).
$rng = $d.Content
$found = $rng.Find.Execute("The creator pattern, where different objects are created")
if (-not $found) {
    throw "Could not locate the 'Low Coupling' paragraph to edit."
}
$rng.Expand(4) | Out-Null

# --- Recover the paragraph's real OOXML (so paraId/rsid/pPr survive) -------
# Range.WordOpenXML returns the whole package; pull just this <w:p>..</w:p>
# back out of it using the unique "Low Coupling" marker text as an anchor
# (NB: String.LastIndexOf(value, startIndex) in this host searches forward,
# not backward like real .NET, so we truncate first and use the 1-arg form).
$wx = $rng.WordOpenXML
$markerIdx = $wx.IndexOf("Low Coupling")
$prefix = $wx.Substring(0, $markerIdx)
$pStart = $prefix.LastIndexOf("<w:p ")
$pEndIdx = $wx.IndexOf("</w:p>", $markerIdx)
$fragment = $wx.Substring($pStart, ($pEndIdx + 6) - $pStart)

# Opening <w:p ...> tag (attributes preserved verbatim).
$pTagEnd = $fragment.IndexOf(">")
$pOpenTag = $fragment.Substring(0, $pTagEnd + 1)
$rest = $fragment.Substring($pTagEnd + 1)

# <w:pPr>...</w:pPr> (preserved verbatim).
$pPrEnd = $rest.IndexOf("</w:pPr>") + 8
$pPr = $rest.Substring(0, $pPrEnd)
$afterPPr = $rest.Substring($pPrEnd)

# First run in the paragraph is the bold "Low Coupling" label; keep as-is.
$firstRunEnd = $afterPPr.IndexOf("</w:r>") + 6
$firstRun = $afterPPr.Substring(0, $firstRunEnd)

# --- Build the replacement body runs ----------------------------------------
$rPr = "<w:rPr><w:sz w:val=`"26`"/><w:szCs w:val=`"26`"/></w:rPr>"

function New-TextRun([string]$text, [bool]$preserveSpace) {
    if ($preserveSpace) {
        return "<w:r>$rPr<w:t xml:space=`"preserve`">$text</w:t></w:r>"
    } else {
        return "<w:r>$rPr<w:t>$text</w:t></w:r>"
    }
}

$spellStart = '<w:proofErr w:type="spellStart"/>'
$spellEnd = '<w:proofErr w:type="spellEnd"/>'

$newBody = ""
$newBody += New-TextRun ": " $true
$newBody += New-TextRun "By using classes such as " $true
$newBody += $spellStart
$newBody += New-TextRun "BidController" $false
$newBody += $spellEnd
$newBody += New-TextRun " and " $true
$newBody += $spellStart
$newBody += New-TextRun "collorPallete" $false
$newBody += $spellEnd
$newBody += New-TextRun " we are lowering the responsibilities of the " $true
$newBody += $spellStart
$newBody += New-TextRun "gameGUI" $false
$newBody += $spellEnd
$newBody += New-TextRun ". This reassignment of responsibilities allows the " $true
$newBody += $spellStart
$newBody += New-TextRun "gameGUI" $false
$newBody += $spellEnd
$newBody += New-TextRun " to be changed without affecting related classes. This is a clear example where the low coupling design pattern was implemented in our design" $true
$newBody += '<w:bookmarkStart w:id="0" w:name="_GoBack"/>'
$newBody += '<w:bookmarkEnd w:id="0"/>'
$newBody += New-TextRun "." $false

# --- Reassemble and write the paragraph back --------------------------------
$newXml = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$fullPTag = $pOpenTag.Substring(0, $pOpenTag.Length - 1) + " " + $newXml + ">"

$newParagraph = $fullPTag + $pPr + $firstRun + $newBody + "</w:p>"

$rng.InsertXML($newParagraph)

# NB: $rng's Start/End do not get re-synced to the newly inserted content by
# this host, so re-Find the paragraph fresh to report/verify what landed.
$check = $d.Content
$check.Find.Execute("By using classes such as") | Out-Null
$check.Expand(4) | Out-Null
Write-Output "Low Coupling paragraph rewritten: $($check.Text)"
